$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.309974789619446
$ws.Range("B1").Value = 1.672612547874451
$ws.Range("C1").Value = 4.089252471923828
$ws.Range("D1").Value = 3.246322631835938
$ws.Range("E1").Value = 1.123023271560669
